$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1877022653721683
$ws.Range("C2").Value = 0.5728155339805825
$ws.Range("J2").Value = 0.006472491909385114
$ws.Range("P2").Value = 0.1326860841423948
$ws.Range("S2").Value = 0.1003236245954693
$ws.Range("B3").Value = 0.01098901098901099
$ws.Range("C3").Value = 0.02197802197802198
$ws.Range("J3").Value = 0.03296703296703297
$ws.Range("P3").Value = 0.7637362637362637
$ws.Range("S3").Value = 0.1703296703296703
$ws.Range("J4").Value = 0.02325581395348837
$ws.Range("P4").Value = 0.7906976744186046
$ws.Range("S4").Value = 0.186046511627907
$ws.Range("B6").Value = 0.05882352941176471
$ws.Range("D6").Value = 0.009049773755656109
$ws.Range("F6").Value = 0.03619909502262444
$ws.Range("J6").Value = 0.2579185520361991
$ws.Range("O6").Value = 0.02714932126696833
$ws.Range("Q6").Value = 0.167420814479638
$ws.Range("R6").Value = 0.04524886877828054
$ws.Range("S6").Value = 0.3981900452488688
$ws.Range("B7").Value = 0.0989010989010989
$ws.Range("D7").Value = 0.02197802197802198
$ws.Range("F7").Value = 0.04945054945054945
$ws.Range("J7").Value = 0.1538461538461539
$ws.Range("O7").Value = 0.03296703296703297
$ws.Range("Q7").Value = 0.1428571428571428
$ws.Range("R7").Value = 0.08791208791208792
$ws.Range("S7").Value = 0.4120879120879121
$ws.Range("B8").Value = 0.09603340292275574
$ws.Range("D8").Value = 0.02713987473903967
$ws.Range("E8").Value = 0.00208768267223382
$ws.Range("F8").Value = 0.06889352818371608
$ws.Range("J8").Value = 0.09812108559498957
$ws.Range("O8").Value = 0.0104384133611691
$ws.Range("Q8").Value = 0.1920668058455115
$ws.Range("R8").Value = 0.08977035490605428
$ws.Range("S8").Value = 0.4154488517745303
$ws.Range("B9").Value = 0.1044776119402985
$ws.Range("D9").Value = 0.01119402985074627
$ws.Range("F9").Value = 0.06716417910447761
$ws.Range("J9").Value = 0.08582089552238806
$ws.Range("O9").Value = 0.01119402985074627
$ws.Range("Q9").Value = 0.2537313432835821
$ws.Range("R9").Value = 0.09328358208955224
$ws.Range("S9").Value = 0.373134328358209
$ws.Range("B10").Value = 0.09985835694050992
$ws.Range("D10").Value = 0.01558073654390935
$ws.Range("E10").Value = 0.00141643059490085
$ws.Range("F10").Value = 0.06586402266288952
$ws.Range("J10").Value = 0.1161473087818697
$ws.Range("O10").Value = 0.009206798866855524
$ws.Range("Q10").Value = 0.2110481586402266
$ws.Range("R10").Value = 0.09702549575070822
$ws.Range("S10").Value = 0.3838526912181303
$ws.Range("F11").Value = 0.003344481605351171
$ws.Range("G11").Value = 0.1705685618729097
$ws.Range("J11").Value = 0.07357859531772576
$ws.Range("K11").Value = 0.2173913043478261
$ws.Range("L11").Value = 0.5183946488294314
$ws.Range("S11").Value = 0.01672240802675585
$ws.Range("G12").Value = 0.7142857142857143
$ws.Range("J12").Value = 0.2111801242236025
$ws.Range("K12").Value = 0.0124223602484472
$ws.Range("L12").Value = 0.04347826086956522
$ws.Range("S12").Value = 0.01863354037267081
$ws.Range("G13").Value = 0.5238095238095238
$ws.Range("J13").Value = 0.4047619047619048
$ws.Range("S13").Value = 0.07142857142857142
$ws.Range("F15").Value = 0.01298701298701299
$ws.Range("H15").Value = 0.1471861471861472
$ws.Range("I15").Value = 0.07792207792207792
$ws.Range("J15").Value = 0.4112554112554113
$ws.Range("K15").Value = 0.03463203463203463
$ws.Range("M15").Value = 0.01298701298701299
$ws.Range("O15").Value = 0.04329004329004329
$ws.Range("S15").Value = 0.2597402597402597
$ws.Range("H16").Value = 0.1707317073170732
$ws.Range("I16").Value = 0.07804878048780488
$ws.Range("J16").Value = 0.375609756097561
$ws.Range("K16").Value = 0.1365853658536585
$ws.Range("M16").Value = 0.02439024390243903
$ws.Range("N16").Value = 0.004878048780487805
$ws.Range("O16").Value = 0.07317073170731707
$ws.Range("S16").Value = 0.1365853658536585
$ws.Range("F17").Value = 0.01356589147286822
$ws.Range("H17").Value = 0.1608527131782946
$ws.Range("I17").Value = 0.1162790697674419
$ws.Range("J17").Value = 0.4186046511627907
$ws.Range("K17").Value = 0.08527131782945736
$ws.Range("M17").Value = 0.01550387596899225
$ws.Range("O17").Value = 0.05232558139534884
$ws.Range("S17").Value = 0.1375968992248062
$ws.Range("F18").Value = 0.0128755364806867
$ws.Range("H18").Value = 0.1416309012875537
$ws.Range("I18").Value = 0.1201716738197425
$ws.Range("J18").Value = 0.48068669527897
$ws.Range("K18").Value = 0.06866952789699571
$ws.Range("M18").Value = 0.02575107296137339
$ws.Range("O18").Value = 0.04721030042918455
$ws.Range("S18").Value = 0.1030042918454936
$ws.Range("F19").Value = 0.01256983240223464
$ws.Range("H19").Value = 0.2074022346368715
$ws.Range("I19").Value = 0.1005586592178771
$ws.Range("J19").Value = 0.3729050279329609
$ws.Range("K19").Value = 0.09567039106145252
$ws.Range("M19").Value = 0.01466480446927374
$ws.Range("O19").Value = 0.07402234636871509
$ws.Range("S19").Value = 0.1222067039106145
